$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: nexial.pollWaitMs value changes from 800 to 0
$ws.Range("B3").Value = "0"

# Row 4: previously nexial.failFast / false -> now nexial.delayBetweenStepsMs / 0
$ws.Range("A4").Value = "nexial.delayBetweenStepsMs"
$ws.Range("B4").Value = "0"

# Rows 5-7: previously held nexial.textDelim/","  nexial.verbose/"false"  nexial.openResult/"true"
# These rows are now emptied out (A cell content cleared, B cell removed entirely)
$ws.Range("A5:A7").ClearContents()
$ws.Range("B5:B7").Clear()

# Update the active selection to A5 as reflected in the sheet view
$ws.Range("A5").Select()
